$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.538.63"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.917.88"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'245.46"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.4838"
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("D8").Value = "'0.2901"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").Value = "'0.06717"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").Value = "'111.21"
$ws.Range("E10").Value = "  +5.35%  "
$ws.Range("D11").Value = "'18.94"
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("D12").Value = "1.918.68"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "'0.07561"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").Value = "'5.290"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "'0.6704"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "'298.43"
$ws.Range("E16").Value = "  +3.98%  "
$ws.Range("D17").Value = "30.530.85"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000007584"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.554"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").Value = "2.166.86"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'0.9999"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'6.427"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("D25").Value = "'9.466"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "'165.16"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").Value = "'20.28"
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("D28").Value = "'2.112"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "'0.1064"
$ws.Range("E29").Value = "  -2.44%  "
$ws.Range("D30").Value = "'1.430"
$ws.Range("E30").Value = "  +5.05%  "
$ws.Range("D31").Value = "'4.147"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").Value = "'4.080"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("D33").Value = "'0.05010"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("D34").Value = "'0.7382"
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("D35").Value = "'1.139"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("D36").Value = "'0.9999"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "'2.727"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'0.02026"
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("D39").Value = "'2.682"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").Value = "'2.016"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("D42").Value = "'0.4446"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "'0.8660"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").Value = "'71.11"
$ws.Range("E44").Value = "  +4.70%  "
$ws.Range("D45").Value = "'5.835"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'48.98"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.216"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").Value = "'9.196"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("D50").Value = "'0.1233"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'34.86"
$ws.Range("E51").Value = "  -1.26%  "
